# Update the date in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-09-11 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-09-12 Tuesday", 2)

# Update the division problems in the table. The table has 20 rows; the
# 5 "data" rows (1, 5, 9, 13, 17) each hold 5 problems in their cells,
# the rest are blank spacer rows.
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("19÷4=", "23÷5=", "96÷8=", "83÷6=", "35÷7=")
    5  = @("66÷5=", "44÷9=", "94÷6=", "82÷6=", "87÷9=")
    9  = @("11÷4=", "38÷9=", "33÷4=", "23÷9=", "12÷9=")
    13 = @("44÷7=", "81÷7=", "74÷7=", "55÷8=", "32÷9=")
    17 = @("37÷7=", "32÷5=", "50÷3=", "38÷5=", "84÷3=")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($row, $c)
        $cell.Range.Text = $vals[$c - 1]
    }
}

Write-Output "done"
